$d = $word.ActiveDocument

# Grab the last paragraph ("public Static int addup(int… numbers){}") and
# append new paragraphs with the new notes after it.
$lastPara = $d.Paragraphs.Last

$lastPara.Range.InsertParagraphAfter()
$r1 = $d.Paragraphs.Last.Range
$r1.InsertAfter("* Classes can also store attributes. ")

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$r2 = $d.Paragraphs.Last.Range
$r2.InsertAfter("* Getters are return statements. They are used to make private attributes public. ")

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$r3 = $d.Paragraphs.Last.Range
$r3.InsertAfter("* Setters are void operations. U can use them to allow user to make changes on app. ")

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$r4 = $d.Paragraphs.Last.Range
$r4.InsertAfter("* this function states the variable in assigned class.  ")
